$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
# A8: "Volume 33   Number  3" -> "...Number  4"  (only the last run, "3", changes to "4")
$ws.Range("A8").Characters(21, 1).Text = "4"

# C9: "Report Covering the Week  1/12/2026  Through  1/18/2026"
#     -> "...1/19/2026  Through  1/25/2026"
$ws.Range("C9").Characters(27, 9).Text = "1/19/2026"
$ws.Range("C9").Characters(47, 9).Text = "1/25/2026"

# --- Crime Complaints data table (rows 14-33) ---
# --- Simple value updates (style unchanged) ---
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 100
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 100
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -55.555555555555
$ws.Range("F16").Value = 16
$ws.Range("H16").Value = -38.461538461538
$ws.Range("I16").Value = 15
$ws.Range("J16").Value = 23
$ws.Range("K16").Value = -34.782608695652
$ws.Range("L16").Value = -25
$ws.Range("M16").Value = -50
$ws.Range("N16").Value = -88.805970149253
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 54
$ws.Range("G17").Value = 61
$ws.Range("H17").Value = -11.475409836065
$ws.Range("I17").Value = 48
$ws.Range("J17").Value = 54
$ws.Range("K17").Value = -11.111111111111
$ws.Range("L17").Value = 29.729729729729
$ws.Range("M17").Value = 242.857142857143
$ws.Range("N17").Value = -23.809523809523
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 15
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = -16.666666666666
$ws.Range("L18").Value = 50
$ws.Range("M18").Value = 7.142857142857
$ws.Range("N18").Value = -79.729729729729
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 87.5
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -2.325581395348
$ws.Range("I19").Value = 39
$ws.Range("J19").Value = 38
$ws.Range("K19").Value = 2.631578947368
$ws.Range("L19").Value = -13.333333333333
$ws.Range("M19").Value = 30
$ws.Range("N19").Value = -41.791044776119
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 12
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 11
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = -15.384615384615
$ws.Range("L20").Value = -47.619047619047
$ws.Range("M20").Value = 22.222222222222
$ws.Range("N20").Value = -88.172043010752
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -5.882352941176
$ws.Range("F21").Value = 147
$ws.Range("G21").Value = 170
$ws.Range("H21").Value = -13.529411764705
$ws.Range("I21").Value = 134
$ws.Range("J21").Value = 149
$ws.Range("K21").Value = -10.067114093959
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 34
$ws.Range("N21").Value = -69.195402298850
$ws.Range("I22").Value = 4
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 100
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = -50
$ws.Range("L23").Value = -50
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -2.941176470588
$ws.Range("F24").Value = 162
$ws.Range("G24").Value = 148
$ws.Range("H24").Value = 9.459459459459
$ws.Range("I24").Value = 143
$ws.Range("J24").Value = 124
$ws.Range("K24").Value = 15.322580645161
$ws.Range("L24").Value = 20.168067226890
$ws.Range("M24").Value = 88.157894736842
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 86
$ws.Range("H25").Value = 24.637681159420
$ws.Range("I25").Value = 78
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = 30
$ws.Range("L25").Value = 16.417910447761
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = -6.25
$ws.Range("F26").Value = 82
$ws.Range("G26").Value = 70
$ws.Range("H26").Value = 17.142857142857
$ws.Range("I26").Value = 76
$ws.Range("J26").Value = 62
$ws.Range("K26").Value = 22.580645161290
$ws.Range("L26").Value = 46.153846153846
$ws.Range("M26").Value = 18.75
$ws.Range("C27").Value = 3
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 9
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 3
$ws.Range("L27").Value = 200
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 5
$ws.Range("J28").Value = 6
$ws.Range("K28").Value = -16.666666666666
$ws.Range("L28").Value = -16.666666666666
$ws.Range("F31").Value = 2
$ws.Range("I31").Value = 2

# --- Cells that change from text ("n/a"/"***.*") to numeric: set NumberFormat to reuse existing style, then Value ---
$ws.Range("N14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N14").Value = -100
$ws.Range("L15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L15").Value = 500
$ws.Range("D33").NumberFormat = "#,##0"
$ws.Range("D33").Value = 1
$ws.Range("E33").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E33").Value = -100
$ws.Range("G33").NumberFormat = "#,##0"
$ws.Range("G33").Value = 1
$ws.Range("H33").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H33").Value = 0
$ws.Range("J33").NumberFormat = "#,##0"
$ws.Range("J33").Value = 1
$ws.Range("K33").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K33").Value = -100
